$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Geographic Distance"
$ws.Range("B8").Value = '<span style="color:#fc8d62">Britto_2016</span>'

$ws.Range("A9").Value = "Temporal Distance"
$ws.Range("B9").Value = '<span style="color:#fc8d62">Britto_2016</span>, <span style="color:#66c2a5">Bajta</span>'

$ws.Range("A10").Value = "Legal Entity"
$ws.Range("B10").Value = '<span style="color:#fc8d62">Britto_2016</span>, <span style="color:#e78ac3">Dashti</span>'
